$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 — shifts existing rows 6:83 down to 7:84
# and extends the used range to A1:R84 (matches the dimension change in
# the diff).
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new data record. Columns
# that are identical to the (now shifted-down) neighbouring rows are
# filled with the same values Excel would have shown for this market /
# product combination; columns called out in the diff get the new
# values.
$ws.Range("A6").Value = 11
$ws.Range("B6").Value = "Vega Monumental Concepción"
$ws.Range("C6").Value = "Bíobío"
$ws.Range("D6").Value = 44761
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 100112012
$ws.Range("G6").Value = "Espinaca"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 45
$ws.Range("K6").Value = 10000
$ws.Range("L6").Value = 11000
$ws.Range("M6").Value = 10333
$ws.Range("N6").Value = "$/cuna 10 kilos"
$ws.Range("O6").Value = "Región Metropolitana"
$ws.Range("P6").Value = 1033
$ws.Range("Q6").Value = 10
$ws.Range("R6").Value = "Hortaliza"
